$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows before row 135, pushing the existing 135:148 rows
# (and their formatting) down to 141:154.
$ws.Range("A135:R140").Insert()

# Common columns shared by every data row in this sheet.
$mercadoId = 4
$mercado = "Feria Lagunitas de Puerto Montt"
$region = "Los Lagos"
$codreg = 10
$categoriaId = 100112027
$categoria = "Melón"
$clasificacion = "Hortaliza"

function Set-Row($r, $fecha, $variedad, $calidad, $volumen, $precioMin, $precioMax, $precioProm, $unidadComerc, $origen, $precioKg, $kgOUnidades) {
    $ws.Cells.Item($r, 1).Value = $mercadoId
    $ws.Cells.Item($r, 2).Value = $mercado
    $ws.Cells.Item($r, 3).Value = $region
    $ws.Cells.Item($r, 4).Value = $fecha
    $ws.Cells.Item($r, 5).Value = $codreg
    $ws.Cells.Item($r, 6).Value = $categoriaId
    $ws.Cells.Item($r, 7).Value = $categoria
    $ws.Cells.Item($r, 8).Value = $variedad
    $ws.Cells.Item($r, 9).Value = $calidad
    $ws.Cells.Item($r, 10).Value = $volumen
    $ws.Cells.Item($r, 11).Value = $precioMin
    $ws.Cells.Item($r, 12).Value = $precioMax
    $ws.Cells.Item($r, 13).Value = $precioProm
    $ws.Cells.Item($r, 14).Value = $unidadComerc
    $ws.Cells.Item($r, 15).Value = $origen
    $ws.Cells.Item($r, 16).Value = $precioKg
    $ws.Cells.Item($r, 17).Value = $kgOUnidades
    $ws.Cells.Item($r, 18).Value = $clasificacion
}

Set-Row 135 44551 "Calameño" "Extra"   300  15000 15000 15000 '$/caja 12 unidades' "Región de O'Higgins" 1250 12
Set-Row 136 44551 "Calameño" "Primera" 1500 1500  1500  1500  '$/unidad'            "Región de O'Higgins" 1500 1
Set-Row 137 44551 "Calameño" "Segunda" 1500 1000  1000  1000  '$/unidad'            "Región de O'Higgins" 1000 1
Set-Row 138 44551 "Tuna"     "Extra"   300  14000 14000 14000 '$/caja 12 unidades' "Región de O'Higgins" 1167 12
Set-Row 139 44551 "Tuna"     "Primera" 1500 1500  1500  1500  '$/unidad'            "Región de O'Higgins" 1500 1
Set-Row 140 44551 "Tuna"     "Segunda" 1500 1000  1000  1000  '$/unidad'            "Región de O'Higgins" 1000 1

# Make sure the D-column cells keep the date number format used elsewhere
# in that column (style index 2 in the original workbook).
$ws.Range("D135:D140").NumberFormat = $ws.Range("D141").NumberFormat

Write-Output "done"
